# Update database: roll the quarterly table one column to the left (drop the
# oldest trailing period) and append the newly published period
# ("12 ماهه منتهی به 1401/12") as the new rightmost column (M).
# Also corrects the amended publish-date footnote for the 1400/12 annual
# report column (row 9), which isn't a pure shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D(4) .. M(13) hold the ten trailing quarterly periods.
$firstCol = 4
$lastCol = 13

# Rows whose 10 period-columns simply shift one column to the left (the
# value that was in column E moves to D, F moves to E, ... and M is freed up
# for the newly appended period). Row 9 (publish dates) is handled
# separately below since it also carries a footnote amendment.
$shiftRows = 8,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27

foreach ($r in $shiftRows) {
    for ($c = $firstCol; $c -lt $lastCol; $c++) {
        $v = $ws.Cells.Item($r, $c + 1).Value2
        $ws.Cells.Item($r, $c).Value2 = $v
    }
}

# Newly appended column M values for each shifted row (the freshly published
# 12-month period ended 1401/12).
$ws.Cells.Item(8, 13).Value2 = "12 ماهه منتهی به 1401/12"
$ws.Cells.Item(11, 13).Value2 = 8941755
$ws.Cells.Item(12, 13).Value2 = -7828334
$ws.Cells.Item(13, 13).Value2 = 1113421
$ws.Cells.Item(14, 13).Value2 = -53912
$ws.Cells.Item(16, 13).Value2 = 4467
$ws.Cells.Item(17, 13).Value2 = 1063975
$ws.Cells.Item(18, 13).Value2 = -7161
$ws.Cells.Item(19, 13).Value2 = 57351
$ws.Cells.Item(20, 13).Value2 = 1114165
$ws.Cells.Item(21, 13).Value2 = -123733
$ws.Cells.Item(22, 13).Value2 = 990432
$ws.Cells.Item(24, 13).Value2 = 990432
$ws.Cells.Item(26, 13).Value2 = 380690

# Row 9 (publish dates) is not a pure shift: the 1400/12 annual report's
# published-date footnote was amended from revision (7) to revision (8) with
# a new date, in addition to the new trailing column being appended. These
# are plain text labels (not real dates), so a leading apostrophe forces
# text entry the same way typing them in Excel would.
$ws.Cells.Item(9, 4).Value2 = "'1400-10-30 (2)"
$ws.Cells.Item(9, 5).Value2 = "'1401-04-27 (10)"
$ws.Cells.Item(9, 6).Value2 = "'1401-04-30 (2)"
$ws.Cells.Item(9, 7).Value2 = "'1401-08-30 (4)"
$ws.Cells.Item(9, 8).Value2 = "'1401-10-29 (2)"
$ws.Cells.Item(9, 9).Value2 = "'1402-02-25 (8)"
$ws.Cells.Item(9, 10).Value2 = "'1401-04-30"
$ws.Cells.Item(9, 11).Value2 = "'1401-08-30 (2)"
$ws.Cells.Item(9, 12).Value2 = "'1401-10-29"
$ws.Cells.Item(9, 13).Value2 = "'1402-02-25"
